$wb = $excel.ActiveWorkbook

# "TestData" sheet: iteration-2 row previously described the "Pca Group
# Profile" case; replace it with the new "Comptroller source/group profile"
# case (same AppYear, same Iteration number).
$wsData = $wb.Worksheets.Item("TestData")

# Set ProfileType (column C) before Description (column B) so the shared
# string table gets the two new strings appended in that order.
$wsData.Range("C3").Value = "COMPTROLLER SOURCE/GROUP PROFILE"
$wsData.Range("B3").Value = "Select COMPTROLLER SOURCE/GROUP PROFILE"

# Row 4 becomes the new iteration-3 placeholder: clear the old empty
# formatted cells in B:D and give A4 the next iteration number.
$wsData.Range("B4:D4").Clear() | Out-Null
$wsData.Range("A4").Value = 3

# Update the recorded selection to reflect the cell last worked on.
$wsData.Activate()
$wsData.Range("D3").Select() | Out-Null
